$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @("Drew Barklay", "Ryan Gustafson", "Alex Shaffer", "Augie Phelps", "Marcus Berger", "Jack Stonis")

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $names[$i]
}

$ws.Range("B2:B7").Select()
